# Adding test cases to Search Module
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# 1) Capture the formatting of row 96 (A:E) BEFORE any edits, so the brand
#    new row 100 can reuse the same "leftover" style pattern it had in the
#    source workbook (A/B/C = s7/s6/s4, D/E = s7).
$ws.Range("A96:E96").Copy()
$ws.Range("A100:E100").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Normalize the D-column style on rows 95-98 (was a stray "fill" style,
#    becomes the plain bordered style used throughout the rest of the sheet).
$ws.Range("D93").Copy()
$ws.Range("D95:D98").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) The "Results" column (E) moves every one of these rows to SKIP - the
#    module's automated run no longer reports individual PASS/FAIL there.
$ws.Range("E93").Value2 = "SKIP"
$ws.Range("E94").Value2 = "SKIP"
$ws.Range("E95").Value2 = "SKIP"
$ws.Range("E96").Value2 = "SKIP"
$ws.Range("E97").Value2 = "SKIP"
$ws.Range("E98").Value2 = "SKIP"
$ws.Range("E99").Value2 = "SKIP"

# 4) New test case TestCase_B99 appended as row 100.
$ws.Range("A100").Value2 = "TestCase_B99"
$ws.Range("C100").Value2 = "Verify that following options get displayed in SORT BY drop down in PEOPLE search results page: a)Relevance b)Registration Date and search results are`nsorted by Relevance by default."
$ws.Range("B100").Value2 = "OPQA-1240 |OPQA-1241"
$ws.Range("D100").Value2 = "Y"
$ws.Range("E100").Value2 = "PASS"
$ws.Rows.Item(100).RowHeight = 45

# 5) Update the sheet view: scroll/selection moved from row 96 area to the
#    D column, with the window's top-left anchored near C88.
$ws.Range("D2:D99").Select()

Write-Host "done"
